$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 39 and 40 (Handajega, block W2) held two duplicate soil-nutrient
# observations. Average the N.per (F) and C.per (G) readings from the two
# rows into row 39 (keep the P.per/H value as-is), then remove the now
# redundant row 40, shifting the remaining observations up by one row.

$f39 = $ws.Cells.Item(39, 6).Value()
$f40 = $ws.Cells.Item(40, 6).Value()
$g39 = $ws.Cells.Item(39, 7).Value()
$g40 = $ws.Cells.Item(40, 7).Value()

$ws.Cells.Item(39, 6).Value = ($f39 + $f40) / 2
$ws.Cells.Item(39, 7).Value = ($g39 + $g40) / 2

# Document the averaging with review comments on the two averaged cells.
$note = "Had two results for same labeled soil. Did an avarage of these two."
$cF = $ws.Range("F39").AddComment($note)
$cG = $ws.Range("G39").AddComment($note)

# Remove the duplicate row; everything below shifts up.
$ws.Rows.Item(40).Delete()

# Re-apply the C:N ratio formula over the (now shorter) data range so it
# stays a single shared formula covering I3:I57.
$ws.Range("I3:I57").Formula = "=G3/F3"

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("L11").Select()
